$wb = $excel.ActiveWorkbook

# 1. Rename the "SID Allocation" sheet to "SID-mapping-bis".
#    Excel automatically updates formula references (e.g. on the
#    "excluded SID" sheet) that point at the old sheet name.
$ws1 = $wb.Worksheets.Item("SID Allocation")
$ws1.Name = "SID-mapping-bis"

# 2. Update a few text values that now reference the plain
#    "/ietf-schc:schc/rule/..." path instead of the
#    "ietf-schc-compound-ack:"-prefixed one, and drop stray
#    trailing spaces that were cleaned up in this revision.
$ws1.Cells.Item(115, 3).Value = "/ietf-schc:schc/rule/bitmap-format"
$ws1.Cells.Item(116, 2).Value = "data"
$ws1.Cells.Item(116, 3).Value = "/ietf-schc:schc/rule/last-bitmap-compression"
$ws1.Cells.Item(241, 2).Value = "identity"

# 3. Restore the view/selection state recorded for this sheet
#    (scrolled so row 228 is at the top, with C246 selected).
$ws1.Range("C246").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 228
